$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: "추천 모델과 LLM - 1" -> "추천 모델과 LLM - 2"
$ws.Range("D7").Value = "추천 모델과 LLM - 2"
$ws.Range("E7").Value = "https://jayhey.github.io/deep%20learning/2023/05/28/reco_llm_2/"

# Row 9: "[대학원] 향후 1년간 SIAI 로드맵" -> "[스타트업] VC업무에 필요한 Finance 관련 지식은 많지 않다?"
$ws.Range("D9").Value = "[스타트업] VC업무에 필요한 Finance 관련 지식은 많지 않다?"
$ws.Range("E9").Value = "https://pdsi.pabii.com/startup-corporate-finance-for-vc/#utm_source=rss&utm_medium=rss&utm_campaign=startup-corporate-finance-for-vc"

# Row 44: "투자 Exit 전략/동향 및 세컨드리 펀드" -> "차량용 반도체 종류와 시장 동향"
$ws.Range("D44").Value = "차량용 반도체 종류와 시장 동향"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/132"

# Row 51: "헷갈리는 DB 용어, DDL, DML, DCL 정리" -> "[vscode] 설정(Settings) 페이지 단축키, ctrl + 콤마"
$ws.Range("D51").Value = "[vscode] 설정(Settings) 페이지 단축키, ctrl + 콤마"
$ws.Range("E51").Value = "https://bskyvision.com/entry/vscode-%EC%84%A4%EC%A0%95Settings-%ED%8E%98%EC%9D%B4%EC%A7%80-%EB%8B%A8%EC%B6%95%ED%82%A4-ctrl-%EC%BD%A4%EB%A7%88"
